$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.748197
$ws.Range("H2").Value = 5.244591
$ws.Range("I2").Value = 0.234333233099101
$ws.Range("J2").Value = 0.234333233099101
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 4.705367
$ws.Range("N2").Value = 14.116101
$ws.Range("O2").Value = 0.08150174587488268
$ws.Range("P2").Value = 0.08150174587488268
$ws.Range("Q2").Value = 8.225908473299
$ws.Range("R2").Value = 74.033176259691
$ws.Range("S2").Value = 0.01909856761408258
$ws.Range("T2").Value = 0.01909856761408257

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.748197
$ws.Range("H3").Value = 5.244591
$ws.Range("I3").Value = 0.234333233099101
$ws.Range("J3").Value = 0.234333233099101
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 39.13297533333333
$ws.Range("N3").Value = 117.398926
$ws.Range("O3").Value = 0.6778229649133395
$ws.Range("P3").Value = 0.6778229649133395
$ws.Range("Q3").Value = 68.41215007880733
$ws.Range("R3").Value = 615.709350709266
$ws.Range("S3").Value = 0.1588364468369614
$ws.Range("T3").Value = 0.1588364468369613

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.748197
$ws.Range("H4").Value = 5.244591
$ws.Range("I4").Value = 0.234333233099101
$ws.Range("J4").Value = 0.234333233099101
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 13.89498533333333
$ws.Range("N4").Value = 41.684956
$ws.Range("O4").Value = 0.2406752892117778
$ws.Range("P4").Value = 0.2406752892117778
$ws.Range("Q4").Value = 24.29117167477733
$ws.Range("R4").Value = 218.620545072996
$ws.Range("S4").Value = 0.05639821864805708
$ws.Range("T4").Value = 0.05639821864805707

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 3.162808
$ws.Range("H5").Value = 9.488423999999998
$ws.Range("I5").Value = 0.4239516623765522
$ws.Range("J5").Value = 0.4239516623765521
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 4.705367
$ws.Range("N5").Value = 14.116101
$ws.Range("O5").Value = 0.08150174587488268
$ws.Range("P5").Value = 0.08150174587488268
$ws.Range("Q5").Value = 14.882172390536
$ws.Range("R5").Value = 133.939551514824
$ws.Range("S5").Value = 0.03455280065024782
$ws.Range("T5").Value = 0.03455280065024781

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 3.162808
$ws.Range("H6").Value = 9.488423999999998
$ws.Range("I6").Value = 0.4239516623765522
$ws.Range("J6").Value = 0.4239516623765521
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 39.13297533333333
$ws.Range("N6").Value = 117.398926
$ws.Range("O6").Value = 0.6778229649133395
$ws.Range("P6").Value = 0.6778229649133395
$ws.Range("Q6").Value = 123.7700874480693
$ws.Range("R6").Value = 1113.930787032624
$ws.Range("S6").Value = 0.2873641727720137
$ws.Range("T6").Value = 0.2873641727720136

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 3.162808
$ws.Range("H7").Value = 9.488423999999998
$ws.Range("I7").Value = 0.4239516623765522
$ws.Range("J7").Value = 0.4239516623765521
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 13.89498533333333
$ws.Range("N7").Value = 41.684956
$ws.Range("O7").Value = 0.2406752892117778
$ws.Range("P7").Value = 0.2406752892117778
$ws.Range("Q7").Value = 43.94717077214933
$ws.Range("R7").Value = 395.5245369493439
$ws.Range("S7").Value = 0.1020346889542907
$ws.Range("T7").Value = 0.1020346889542906

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 2.549298333333333
$ws.Range("H8").Value = 7.647895
$ws.Range("I8").Value = 0.3417151045243469
$ws.Range("J8").Value = 0.3417151045243469
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 4.705367
$ws.Range("N8").Value = 14.116101
$ws.Range("O8").Value = 0.08150174587488268
$ws.Range("P8").Value = 0.08150174587488268
$ws.Range("Q8").Value = 11.99538425082167
$ws.Range("R8").Value = 107.958458257395
$ws.Range("S8").Value = 0.02785037761055229
$ws.Range("T8").Value = 0.02785037761055229

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 2.549298333333333
$ws.Range("H9").Value = 7.647895
$ws.Range("I9").Value = 0.3417151045243469
$ws.Range("J9").Value = 0.3417151045243469
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 39.13297533333333
$ws.Range("N9").Value = 117.398926
$ws.Range("O9").Value = 0.6778229649133395
$ws.Range("P9").Value = 0.6778229649133395
$ws.Range("Q9").Value = 99.76162879564112
$ws.Range("R9").Value = 897.85465916077
$ws.Range("S9").Value = 0.2316223453043645
$ws.Range("T9").Value = 0.2316223453043645

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.549298333333333
$ws.Range("H10").Value = 7.647895
$ws.Range("I10").Value = 0.3417151045243469
$ws.Range("J10").Value = 0.3417151045243469
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 13.89498533333333
$ws.Range("N10").Value = 41.684956
$ws.Range("O10").Value = 0.2406752892117778
$ws.Range("P10").Value = 0.2406752892117778
$ws.Range("Q10").Value = 35.42246295195778
$ws.Range("R10").Value = 318.80216656762
$ws.Range("S10").Value = 0.08224238160943007
$ws.Range("T10").Value = 0.08224238160943007
